# Auto-generated Excel COM-interop script
# Updates search-term keyword list in column A (rows 1-100)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "recovery compression pants",
    "knee volleyball",
    "leggings black",
    "mens volleyball pads",
    "knee arthritis compression",
    "youth small baseball pants",
    "mens black pants",
    "knee high leggings",
    "knee pads arthritis",
    "basketball knee pads nike pro",
    "knee compression tights",
    "knee pad sleeve youth",
    "thermal compression pants men winter",
    "honey comb knee pad",
    "combat pants with knee pads",
    "deportivo para hombres",
    "compression spats",
    "red mens compression pants",
    "green basketball knee pads",
    "red basketball leggings",
    "underarmour mens tights",
    "three quarter compression pants men",
    "base layer pants men cold weather",
    "tesla base layer pants",
    "winter gym leggings",
    "thick spandex pants men",
    "knee pad leggings for basketball",
    "youth knee pads for basketball",
    "girls basketball knee pad leggings",
    "knee pads compression",
    "knee pad tights basketball",
    "weightlifting knee pads",
    "compression leggings youth",
    "hex pad knee",
    "knee pad soccer",
    "knee pad for running",
    "basketballs black",
    "boys soccer tights",
    "compression men leggings",
    "youth compression knee pads",
    "leggings with knee",
    "running compression pants men",
    "knee protector volleyball",
    "volleyball knee pads for men",
    "knee guard basketball",
    "black baseball pants boys",
    "compression pants boys",
    "knee pad adult",
    "little boy leggings",
    "capri legging pack",
    "men knee pads",
    "basketball compression",
    "boys compression tight",
    "baseball pants knee high",
    "small volleyball knee pads",
    "knee pads volleyball",
    "youth volleyball",
    "knee protector football",
    "knee pads protector",
    "snowboarding hip pads",
    "knee pad running",
    "men pads",
    "cycling compression tights",
    "capri leggings pack",
    "pants soccer",
    "running compression leggings",
    "boys black baseball pants",
    "leggings for boys",
    "youth girls compression pants",
    "athletic pants soccer",
    "knee guards basketball",
    "compression pants for hockey",
    "wrestling gear for men",
    "adult medium baseball pants",
    "knee pads for arthritis",
    "baseball boy pants",
    "knee protector pad",
    "knee protector soccer",
    "men sports compression",
    "basketball pants boys",
    "baseball pants",
    "knee protector pads",
    "knee pad for sports",
    "sports knee protectors",
    "youth knee pads volleyball",
    "bjj pants",
    "knee pads cycling",
    "boys youth basketball",
    "mens athletic pants tall",
    "legging pack",
    "girl volleyball knee pads",
    "running tights",
    "youth baseball pants",
    "boys sports pants",
    "compression for men",
    "basketballs youth size",
    "knee pads xl",
    "compression knee men",
    "knee guards for adults",
    "football knee pad"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
